$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("B4").Value = 741000000.0
$ws.Range("C4").Value = 867000000.0
$ws.Range("D4").Value = 608000000.0
$ws.Range("E4").Value = 626000000.0
$ws.Range("F4").Value = 777000000.0

# Row 13 - Accounts Payable
$ws.Range("B13").Value = 863000000.0
$ws.Range("C13").Value = 921000000.0
$ws.Range("D13").Value = 492000000.0
$ws.Range("E13").Value = 702000000.0
$ws.Range("F13").Value = 759000000.0

# Row 21 - Long Term Tax Liability (Deferred)
$ws.Range("B21").Value = 206000000.0
$ws.Range("C21").Value = 200000000.0
$ws.Range("D21").Value = 213000000.0
$ws.Range("E21").Value = 214000000.0
$ws.Range("F21").Value = 178000000.0
